$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1648366424899244
$ws.Range("C2").Value = 0.6565107402614302
$ws.Range("D2").Value = 0.6208915960224134
$ws.Range("E2").Value = 0.7879667480436046
$ws.Range("F2").Value = 0.7996194615330313

$ws.Range("B3").Value = 0.06925299554247144
$ws.Range("C3").Value = 0.6317528524474937
$ws.Range("D3").Value = 0.5154239904830427
$ws.Range("E3").Value = 0.7179303521115699
$ws.Range("F3").Value = 0.7437609478291553

$ws.Range("B4").Value = 0.1042866910554764
$ws.Range("C4").Value = 0.5871660518276492
$ws.Range("D4").Value = 0.4576314133407699
$ws.Range("E4").Value = 0.6764845994852875
$ws.Range("F4").Value = 0.6981187964102095

$ws.Range("B5").Value = 0.1411991018086509
$ws.Range("C5").Value = 0.5954485264947712
$ws.Range("D5").Value = 0.500182751991249
$ws.Range("E5").Value = 0.7072359945529136
$ws.Range("F5").Value = 0.7268219329406943

$ws.Range("B6").Value = 0.3268981278322622
$ws.Range("C6").Value = 0.6566970789475719
$ws.Range("D6").Value = 0.6666380539392567
$ws.Range("E6").Value = 0.8164790590941428
$ws.Range("F6").Value = 0.7886526259380041

$ws.Range("B7").Value = 0.2872108118487342
$ws.Range("C7").Value = 0.7140200429738547
$ws.Range("D7").Value = 0.7624163934427893
$ws.Range("E7").Value = 0.8731645855408872
$ws.Range("F7").Value = 0.8745954126766147

$ws.Range("B8").Value = 0.2887732378318851
$ws.Range("C8").Value = 0.6550352748098781
$ws.Range("D8").Value = 0.6159387406405581
$ws.Range("E8").Value = 0.7848176480180336
$ws.Range("F8").Value = 0.7994113517477577

$ws.Range("B9").Value = -0.05480200069631858
$ws.Range("C9").Value = 0.1889080360497501
$ws.Range("D9").Value = 0.05129867587295883
$ws.Range("E9").Value = 0.2264921099574085
$ws.Range("F9").Value = 0.2691526052055958

$ws.Range("B10").Value = -0.1770340777201527
$ws.Range("C10").Value = 0.1770340777201527
$ws.Range("D10").Value = 0.03134106467422507
$ws.Range("E10").Value = 0.1770340777201527
